$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generator at Herat (GB 980845): meter serial number is now known (numeric)
# and KVA Rating (previously blank / EOR null) is now populated.
$ws.Range("A2").Value = 980845
$ws.Range("G2").Value = 150
